$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player/position/team table in the desired row order (row 2 .. row 18)
$data = @(
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("OG Anunoby", "SF,PF", "New York Knicks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
